$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Stash the formatting of the merged "Batch Size: 16" label cell (F21) so we
# can re-apply it later once the merged range grows -- this lets us avoid
# Excel's automatic border redistribution that normally happens when merging
# already-bordered cells.
# ---------------------------------------------------------------------------
$ws.Range("F21").Copy()
$ws.Range("Z100").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Row 17 ("Test" / "Default Params" under the DCNN Experiments table):
# used to be highlighted with the "Good" cell style -- clear it back to Normal.
# ---------------------------------------------------------------------------
$ws.Range("A17:E17").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 18 ("Test" / "Extra Epochs" under the DCNN Experiments table):
# becomes the new highlighted ("Good") row, and its Avg Accuracy (E18) is
# now populated with a value (matches the batch size run -> 1).
# ---------------------------------------------------------------------------
$ws.Range("A18:E18").Style = "Good"
$ws.Range("E18").Value = 1

# ---------------------------------------------------------------------------
# Row 22 ("Test" / "Default Params" under the DRNN Experiments table):
# mirrors row 17 -- clear the "Good" highlight back to Normal.
# ---------------------------------------------------------------------------
$ws.Range("A22:E22").Style = "Normal"

# ---------------------------------------------------------------------------
# Row 23 is a brand new row ("Test" / "Extra Epochs") mirroring row 18,
# added under the DRNN Experiments table so the two tables match again.
# ---------------------------------------------------------------------------
$ws.Range("A23").Value = "Test"
$ws.Range("B23").Value = "Extra Epochs"
$ws.Range("C23").Value = 5
$ws.Range("D23").Value = 50
$ws.Range("E23").Value = 0.81034482853679801
$ws.Range("A23:E23").Style = "Good"

# ---------------------------------------------------------------------------
# Grow the "Batch Size: 16" merged label from F21:F22 to F21:F23, then
# restore its original formatting (border / fill / alignment) uniformly
# across the whole merged range.
# ---------------------------------------------------------------------------
$ws.Range("F21:F22").UnMerge()
$ws.Range("F21:F23").Style = "Normal"
$ws.Range("F21:F23").Merge()
$ws.Range("Z100").Copy()
$ws.Range("F21:F23").PasteSpecial(-4122)
$ws.Range("Z100").Clear()

# ---------------------------------------------------------------------------
# Resize the DRNN Experiments table (Table134 / table3.xml) so it now
# includes the new row 23 in its range and autofilter.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $ws.ListObjects.Count; $i++) {
    $lo = $ws.ListObjects.Item($i)
    if ($lo.Name -eq "Table134") {
        $lo.Resize($ws.Range("A21:E23"))
    }
}

# Update the selection to match where the edit was made.
$ws.Activate()
$ws.Range("F21:F23").Select()
